# Proteome_size.xlsx — spelling fixes to condition labels + cursor move.
#
# "Piruvate" -> "Piruvato", "Xylosa" -> "Xilosa", "Mannosa" -> "Manosa"
# (the shared-string table re-packs itself on save, dropping the now-unused
# old spellings and appending the corrected ones, which also reshuffles the
# <v> indices referenced by the header row — handled automatically by the
# engine, we just need to set the corrected text on the right cells).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("I1").Value = "Piruvato"
$ws.Range("S1").Value = "Xilosa"
$ws.Range("T1").Value = "Manosa"

# Leave the selection where the author last left it when saving.
$ws.Range("T2").Select()
